$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters for shared-string table allocation:
# 1) C3 gets the brand-new "all census tracts in Missouri" text (reuses the freed slot
#    left by overwriting C2 below -> becomes shared string index 4).
# 2) B3 gets the new record id "R12459585" (-> shared string index 5).
# 3) C2 is updated in place to the reworded Illinois+Missouri description
#    (-> appended as shared string index 6, since the old slot was already reused by C3).

$ws.Range("C3").Value = "SocialExplorer.com data for all census tracts in Missouri.  ACS 2018 (5-year estimates) for tables A02001, A01001, A04001, A14028, A12001, A12002, B17008, A14001, A10016, A13003A, A13003B, A13003C, A09001, B09001. Dollar values inflation adjusted to 2018 dollars."
$ws.Range("B3").Value = "R12459585"
$ws.Range("C2").Value = "SocialExplorer.com census tract level data for Illinois and Missouri.  ACS 2018 (5-year estimates) for tables A00001, A00002, A00003, A02001, A01001, A02001, A04001, A12001, B17008, A17009."

# New row's date (16-Feb-2020), matching the same date-style formatting as A2
$ws.Range("A3").Value = 43877
$ws.Range("A3").NumberFormat = "dd\-mmm\-yyyy"
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A3").WrapText = $true

# Match the wrap/alignment formatting used by B2:C2 for the new B3:C3 cells
$ws.Range("B3:C3").HorizontalAlignment = -4131
$ws.Range("B3:C3").VerticalAlignment = -4160
$ws.Range("B3:C3").WrapText = $true

$ws.Rows.Item(3).RowHeight = 45

# Selection moves on to the next empty row, like it did for row 2 -> row 3 before
$ws.Range("A4").Select()
